$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 864
$ws1.Range("F5").Value = 1178
$ws1.Range("F6").Value = 59
$ws1.Range("F7").Value = 4278
$ws1.Range("F10").Value = 2459
$ws1.Range("F15").Value = 654
$ws1.Range("F17").Value = 106
$ws1.Range("F18").Value = 312
$ws1.Range("F19").Value = 24
$ws1.Range("F22").Value = 12
$ws1.Range("F23").Value = 459
$ws1.Range("F26").Value = 515
$ws1.Range("F27").Value = 683
$ws1.Range("F30").Value = 385
$ws1.Range("F33").Value = 947
$ws1.Range("F34").Value = 74
$ws1.Range("F36").Value = 1033
$ws1.Range("F37").Value = 1997
$ws1.Range("F38").Value = 245
$ws1.Range("F39").Value = 6
$ws1.Range("F40").Value = 529
$ws1.Range("F43").Value = 622
$ws1.Range("F44").Value = 1281
$ws1.Range("F45").Value = 72
$ws1.Range("F47").Value = 419

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 864
$ws4.Range("F3").Value = 1178
$ws4.Range("F5").Value = 59
$ws4.Range("F6").Value = 4278
$ws4.Range("F8").Value = 2459
$ws4.Range("F12").Value = 654
$ws4.Range("F14").Value = 106
$ws4.Range("F15").Value = 312
$ws4.Range("F16").Value = 24
$ws4.Range("F19").Value = 459
$ws4.Range("F22").Value = 515
$ws4.Range("F23").Value = 683
$ws4.Range("F29").Value = 385
$ws4.Range("F31").Value = 947
$ws4.Range("F32").Value = 74
$ws4.Range("F35").Value = 1033
$ws4.Range("F36").Value = 1997
$ws4.Range("F37").Value = 245
$ws4.Range("F41").Value = 529
$ws4.Range("F44").Value = 622
$ws4.Range("F45").Value = 1281
$ws4.Range("F46").Value = 72
$ws4.Range("F47").Value = 419
